$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns D, M, N, O, P, R, S between row pairs
# (2,6), (3,7), (4,8), (5,9) -- the "La Ligua" rows and the
# "Provincia de Limarí" rows traded their date/volume/price/origin data.
$cols = @("D", "M", "N", "O", "P", "R", "S")
$pairs = @(@(2, 6), @(3, 7), @(4, 8), @(5, 9))

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
